$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.024.91"
$ws.Range("E2").Value = "'  +0.79%  "
$ws.Range("D3").Value = "'2.245.24"
$ws.Range("E3").Value = "'  +2.32%  "
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("D5").Value = "'271.60"
$ws.Range("E5").Value = "'  +4.72%  "
$ws.Range("D6").Value = "'93.66"
$ws.Range("E6").Value = "'  +14.14%  "
$ws.Range("D7").Value = "'0.623"
$ws.Range("E7").Value = "'  -0.07%  "
$ws.Range("D9").Value = "'0.633"
$ws.Range("E9").Value = "'  +6.47%  "
$ws.Range("D10").Value = "'46.07"
$ws.Range("E10").Value = "'  +6.10%  "
$ws.Range("D11").Value = "'0.0971"
$ws.Range("E11").Value = "'  +5.92%  "
$ws.Range("D12").Value = "'8.30"
$ws.Range("E12").Value = "'  +19.17%  "
$ws.Range("E13").Value = "'  +1.77%  "
$ws.Range("B14").Value = "'Chainlink"
$ws.Range("C14").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'15.25"
$ws.Range("E14").Value = "'  +7.38%  "
$ws.Range("B15").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "'2.584.10"
$ws.Range("E15").Value = "'  +2.56%  "
$ws.Range("E16").Value = "'  +5.13%  "
$ws.Range("D17").Value = "'2.249.18"
$ws.Range("E17").Value = "'  +3.02%  "
$ws.Range("D18").Value = "'43.971.16"
$ws.Range("E18").Value = "'  +0.83%  "
$ws.Range("E19").Value = "'  +2.71%  "
$ws.Range("D20").Value = "'6.17"
$ws.Range("E20").Value = "'  +4.31%  "
$ws.Range("D21").Value = "'70.89"
$ws.Range("E21").Value = "'  +1.23%  "
$ws.Range("E22").Value = "'  -1.94%  "
$ws.Range("D23").Value = "'235.31"
$ws.Range("E23").Value = "'  +2.06%  "
$ws.Range("D24").Value = "'9.14"
$ws.Range("E24").Value = "'  +3.55%  "
$ws.Range("E25").Value = "'  -0.05%  "
$ws.Range("D26").Value = "'11.43"
$ws.Range("E26").Value = "'  +6.94%  "
$ws.Range("E27").Value = "'  +12.61%  "
$ws.Range("D28").Value = "'3.57"
$ws.Range("E28").Value = "'  +6.19%  "
$ws.Range("D29").Value = "'41.15"
$ws.Range("E29").Value = "'  -2.56%  "
$ws.Range("E30").Value = "'  +3.11%  "
$ws.Range("D31").Value = "'172.63"
$ws.Range("E31").Value = "'  -0.43%  "
$ws.Range("D32").Value = "'0.0917"
$ws.Range("E32").Value = "'  +5.43%  "
$ws.Range("D33").Value = "'21.03"
$ws.Range("E33").Value = "'  +3.00%  "
$ws.Range("D34").Value = "'5.49"
$ws.Range("E34").Value = "'  +3.65%  "
$ws.Range("D35").Value = "'0.115"
$ws.Range("E35").Value = "'  +0.65%  "
$ws.Range("E36").Value = "'  +1.85%  "
$ws.Range("D37").Value = "'0.0355"
$ws.Range("E37").Value = "'  +0.96%  "
$ws.Range("D38").Value = "'4.34"
$ws.Range("E38").Value = "'  -3.09%  "
$ws.Range("D39").Value = "'3.60"
$ws.Range("E39").Value = "'  +26.25%  "
$ws.Range("E40").Value = "'  -1.03%  "
$ws.Range("D41").Value = "'0.226"
$ws.Range("E41").Value = "'  +13.79%  "
$ws.Range("E42").Value = "'  +4.07%  "
$ws.Range("E43").Value = "'  +2.12%  "
$ws.Range("D44").Value = "'5.42"
$ws.Range("E44").Value = "'  -0.46%  "
$ws.Range("D45").Value = "'0.0999"
$ws.Range("E45").Value = "'  +1.90%  "
$ws.Range("B46").Value = "'Aave"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'101.09"
$ws.Range("E46").Value = "'  +0.18%  "
$ws.Range("B47").Value = "'FraxShare"
$ws.Range("C47").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'8.40"
$ws.Range("E47").Value = "'  +2.48%  "
$ws.Range("E48").Value = "'  +4.69%  "
$ws.Range("E49").Value = "'  +2.64%  "
$ws.Range("D50").Value = "'0.447"
$ws.Range("E50").Value = "'  +1.88%  "
$ws.Range("D51").Value = "'2.466.89"
$ws.Range("E51").Value = "'  +2.41%  "
